$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb2"
$ws.Range("C2").Value = "Ephb6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 36.899643
$ws.Range("H2").Value = 110.698929
$ws.Range("I2").Value = 0.7238945645409351
$ws.Range("J2").Value = 0.7238945645409351
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2721246666666666
$ws.Range("N2").Value = 0.8163739999999999
$ws.Range("O2").Value = 0.1154759865526449
$ws.Range("P2").Value = 0.1154759865526449
$ws.Range("Q2").Value = 10.041303051494
$ws.Range("R2").Value = 90.37172746344599
$ws.Range("S2").Value = 0.08359243900046173
$ws.Range("T2").Value = 0.08359243900046173

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb2"
$ws.Range("C3").Value = "Ephb6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 36.899643
$ws.Range("H3").Value = 110.698929
$ws.Range("I3").Value = 0.7238945645409351
$ws.Range("J3").Value = 0.7238945645409351
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8319233333333332
$ws.Range("N3").Value = 2.49577
$ws.Range("O3").Value = 0.3530263126440755
$ws.Range("P3").Value = 0.3530263126440755
$ws.Range("Q3").Value = 30.69767400336999
$ws.Range("R3").Value = 276.2790660303299
$ws.Range("S3").Value = 0.2555538288629751
$ws.Range("T3").Value = 0.2555538288629751

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb2"
$ws.Range("C4").Value = "Ephb6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 36.899643
$ws.Range("H4").Value = 110.698929
$ws.Range("I4").Value = 0.7238945645409351
$ws.Range("J4").Value = 0.7238945645409351
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.252499666666667
$ws.Range("N4").Value = 3.757499
$ws.Range("O4").Value = 0.5314977008032796
$ws.Range("P4").Value = 0.5314977008032796
$ws.Range("Q4").Value = 46.216790557619
$ws.Range("R4").Value = 415.951115018571
$ws.Range("S4").Value = 0.3847482966774983
$ws.Range("T4").Value = 0.3847482966774983

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb2"
$ws.Range("C5").Value = "Ephb6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.374819
$ws.Range("H5").Value = 10.124457
$ws.Range("I5").Value = 0.0662069584361419
$ws.Range("J5").Value = 0.0662069584361419
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2721246666666666
$ws.Range("N5").Value = 0.8163739999999999
$ws.Range("O5").Value = 0.1154759865526449
$ws.Range("P5").Value = 0.1154759865526449
$ws.Range("Q5").Value = 0.9183714954353333
$ws.Range("R5").Value = 8.265343458917998
$ws.Range("S5").Value = 0.007645313842063439
$ws.Range("T5").Value = 0.007645313842063439

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb2"
$ws.Range("C6").Value = "Ephb6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.374819
$ws.Range("H6").Value = 10.124457
$ws.Range("I6").Value = 0.0662069584361419
$ws.Range("J6").Value = 0.0662069584361419
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8319233333333332
$ws.Range("N6").Value = 2.49577
$ws.Range("O6").Value = 0.3530263126440755
$ws.Range("P6").Value = 0.3530263126440755
$ws.Range("Q6").Value = 2.807590671876667
$ws.Range("R6").Value = 25.26831604689
$ws.Range("S6").Value = 0.02337279840809074
$ws.Range("T6").Value = 0.02337279840809074

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb2"
$ws.Range("C7").Value = "Ephb6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.374819
$ws.Range("H7").Value = 10.124457
$ws.Range("I7").Value = 0.0662069584361419
$ws.Range("J7").Value = 0.0662069584361419
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.252499666666667
$ws.Range("N7").Value = 3.757499
$ws.Range("O7").Value = 0.5314977008032796
$ws.Range("P7").Value = 0.5314977008032796
$ws.Range("Q7").Value = 4.226959672560334
$ws.Range("R7").Value = 38.042637053043
$ws.Range("S7").Value = 0.03518884618598771
$ws.Range("T7").Value = 0.03518884618598771

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efnb2"
$ws.Range("C8").Value = "Ephb6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.699319
$ws.Range("H8").Value = 32.097957
$ws.Range("I8").Value = 0.2098984770229228
$ws.Range("J8").Value = 0.2098984770229228
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2721246666666666
$ws.Range("N8").Value = 0.8163739999999999
$ws.Range("O8").Value = 0.1154759865526449
$ws.Range("P8").Value = 0.1154759865526449
$ws.Range("Q8").Value = 2.911548616435333
$ws.Range("R8").Value = 26.203937547918
$ws.Range("S8").Value = 0.02423823371011967
$ws.Range("T8").Value = 0.02423823371011967

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efnb2"
$ws.Range("C9").Value = "Ephb6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.699319
$ws.Range("H9").Value = 32.097957
$ws.Range("I9").Value = 0.2098984770229228
$ws.Range("J9").Value = 0.2098984770229228
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8319233333333332
$ws.Range("N9").Value = 2.49577
$ws.Range("O9").Value = 0.3530263126440755
$ws.Range("P9").Value = 0.3530263126440755
$ws.Range("Q9").Value = 8.901013126876666
$ws.Range("R9").Value = 80.10911814189
$ws.Range("S9").Value = 0.07409968537300965
$ws.Range("T9").Value = 0.07409968537300965

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efnb2"
$ws.Range("C10").Value = "Ephb6"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.699319
$ws.Range("H10").Value = 32.097957
$ws.Range("I10").Value = 0.2098984770229228
$ws.Range("J10").Value = 0.2098984770229228
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.252499666666667
$ws.Range("N10").Value = 3.757499
$ws.Range("O10").Value = 0.5314977008032796
$ws.Range("P10").Value = 0.5314977008032796
$ws.Range("Q10").Value = 13.40089348106034
$ws.Range("R10").Value = 120.608041329543
$ws.Range("S10").Value = 0.1115605579397935
$ws.Range("T10").Value = 0.1115605579397935

